$d = $word.ActiveDocument

$newText = "Dates à utiliser pour la Campagne Perseus: 16-25 janvier, 7-16 novembre, 6-15 décembre"

$xmlTemplate = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>__TEXT__</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertXml = $xmlTemplate.Replace("__TEXT__", $newText)

# Repeatedly locate the next paragraph that still holds the old
# "Dates a utiliser pour la Campagne ..." text and rebuild its run content
# as a single clean run (no rPr), leaving the paragraph's own pPr untouched.
# Re-scanning after every rewrite (instead of trusting cached offsets) keeps
# this correct regardless of how much earlier edits shift later offsets.
$safety = 0
while ($safety -lt 50) {
    $safety = $safety + 1
    $found = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "*Dates*utiliser*pour*la*Campagne*octobre*") {
            $found = $p
            break
        }
    }
    if ($found -eq $null) {
        break
    }

    $r = $found.Range
    # Exclude the trailing paragraph mark from the replace range.
    $r.End = $r.End - 1
    # Leave one placeholder character so the range never collapses to zero
    # length (an empty range causes InsertXML to replace the whole <w:p>,
    # wiping out pPr along with it).
    $r.Text = "X"
    $insertionPoint = $d.Range($r.Start, $r.Start)
    $insertionPoint.InsertXML($insertXml)
    # Remove the placeholder character now sitting right after our inserted text.
    $placeholder = $d.Range($r.Start + $newText.Length, $r.Start + $newText.Length + 1)
    $placeholder.Delete()
}
